# RoleBasedAccessManagementData.xlsx - "Extent report code" update
$wb = $excel.ActiveWorkbook

# --- Sheet "Create" ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("B1").Value = "ErrorMsg1"
$wsCreate.Range("A2").Value = "PalakAdmin"
$wsCreate.Range("B2").Value = "Enter Role Name"

# --- Sheet "Edit" ---
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "'PalakAdmin"
$wsEdit.Range("B2").Value = "PalakAdmin Updated"
$wsEdit.Range("C2").Value = "Modified"

# --- Sheet "Delete" ---
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("B2").Value = "deleted"
$wsDelete.Range("A2").Value = "Empire"

# --- Selections / active sheet ---
[void]$wsCreate.Range("B1").Select()
[void]$wsEdit.Range("A1:C2").Select()
[void]$wsDelete.Activate()
[void]$wsDelete.Range("A2").Select()
